# Updated return statistics + new scenario generation method
# Apply new forecast values for column B (MSTL) rows 2-25

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 149.0171813964844
    3  = 147.2691802978516
    4  = 147.1151580810547
    5  = 147.4776000976562
    6  = 153.7435607910156
    7  = 142.9864807128906
    8  = 148.9331207275391
    9  = 147.4106903076172
    10 = 156.2419586181641
    11 = 172.0967864990234
    12 = 178.50146484375
    13 = 222.7992553710938
    14 = 233.2488250732422
    15 = 235.0869445800781
    16 = 208.1131591796875
    17 = 183.5649108886719
    18 = 159.7280883789062
    19 = 150.7450408935547
    20 = 139.3758697509766
    21 = 118.7402420043945
    22 = 116.2283554077148
    23 = 101.978271484375
    24 = 101.1797790527344
    25 = 85.94877624511719
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
